$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert 3 new rows at 314-316 (push blank placeholder rows down) ---
$ws.Range("A314:A316").EntireRow.Insert()

# --- Step 2: remove the now-redundant trailing blank rows (464-467) to restore row count ---
$ws.Range("A464:A467").EntireRow.Delete()

# --- Step 3: update modified data cells in rows 298-313 ---
# Row 298
$ws.Cells.Item(298, 10).Value = 1791
$ws.Cells.Item(298, 11).Value = 0

# Row 300
$ws.Cells.Item(300, 10).Value = 720
$ws.Cells.Item(300, 11).Value = 0

# Row 302
$ws.Cells.Item(302, 10).Value = 977
$ws.Cells.Item(302, 11).Value = 0

# Row 303
$ws.Cells.Item(303, 10).Value = 723
$ws.Cells.Item(303, 11).Value = 0

# Row 304
$ws.Cells.Item(304, 10).Value = 1083
$ws.Cells.Item(304, 11).Value = 0

# Row 305
$ws.Cells.Item(305, 10).Value = 2256
$ws.Cells.Item(305, 11).Value = 0

# Row 307
$ws.Cells.Item(307, 10).Value = 743
$ws.Cells.Item(307, 11).Value = 0
$ws.Cells.Item(307, 15).Value = "T_1"

# Row 309
$ws.Cells.Item(309, 4).Value = 31
$ws.Cells.Item(309, 7).Value = "ifop"
$ws.Cells.Item(309, 8).Value = "included"
$ws.Cells.Item(309, 9).Value = "regular"
$ws.Cells.Item(309, 10).Value = 1000
$ws.Cells.Item(309, 12).Value = 1
$ws.Cells.Item(309, 13).Value = 0.5
$ws.Cells.Item(309, 15).Value = 0.5
$ws.Cells.Item(309, 17).Value = 3.5
$ws.Cells.Item(309, 20).Value = 2
$ws.Cells.Item(309, 22).Value = 27
$ws.Cells.Item(309, 26).Value = 2.5
$ws.Cells.Item(309, 28).Value = 22
$ws.Cells.Item(309, 29).Value = 10.5

# Row 310
$ws.Cells.Item(310, 4).Value = 30
$ws.Cells.Item(310, 6).Value = 2
$ws.Cells.Item(310, 7).Value = "ipsos"
$ws.Cells.Item(310, 8).Value = "excluded"
$ws.Cells.Item(310, 10).Value = 1066
$ws.Cells.Item(310, 11).Value = 0
$ws.Cells.Item(310, 12).Value = 1
$ws.Cells.Item(310, 14).Value = 1.5
$ws.Cells.Item(310, 16).Value = 15.5
$ws.Cells.Item(310, 17).Value = 3
$ws.Cells.Item(310, 19).Value = 6
$ws.Cells.Item(310, 20).Value = 2
$ws.Cells.Item(310, 22).Value = 26
$ws.Cells.Item(310, 26).Value = 2.5
$ws.Cells.Item(310, 27).Value = 1.5
$ws.Cells.Item(310, 28).Value = 21

# Row 311
$ws.Cells.Item(311, 6).Value = 2
$ws.Cells.Item(311, 7).Value = "elabe"
$ws.Cells.Item(311, 8).Value = "partially"
$ws.Cells.Item(311, 10).Value = 999
$ws.Cells.Item(311, 11).Value = 0
$ws.Cells.Item(311, 14).Value = 1.5
$ws.Cells.Item(311, 19).Value = 4.5
$ws.Cells.Item(311, 20).Value = 1.5
$ws.Cells.Item(311, 22).Value = 28.5
$ws.Cells.Item(311, 23).Value = 8.5
$ws.Cells.Item(311, 27).Value = 2.5
$ws.Cells.Item(311, 29).Value = 9.5

# Row 312
$ws.Cells.Item(312, 3).Value = 4
$ws.Cells.Item(312, 4).Value = 2
$ws.Cells.Item(312, 6).Value = 3
$ws.Cells.Item(312, 7).Value = "opinionway"
$ws.Cells.Item(312, 8).Value = "partially"
$ws.Cells.Item(312, 9).Value = "regular"
$ws.Cells.Item(312, 10).Value = 700
$ws.Cells.Item(312, 11).Value = 1
$ws.Cells.Item(312, 13).Value = 1
$ws.Cells.Item(312, 14).Value = 1
$ws.Cells.Item(312, 15).Value = "T_1"
$ws.Cells.Item(312, 16).Value = 14
$ws.Cells.Item(312, 17).Value = 4
$ws.Cells.Item(312, 22).Value = 27
$ws.Cells.Item(312, 23).Value = 10
$ws.Cells.Item(312, 26).Value = 3
$ws.Cells.Item(312, 27).Value = 2
$ws.Cells.Item(312, 28).Value = 22
$ws.Cells.Item(312, 29).Value = 9

# Row 313
$ws.Cells.Item(313, 6).Value = 4
$ws.Cells.Item(313, 7).Value = "opinionway"
$ws.Cells.Item(313, 9).Value = "rolling"
$ws.Cells.Item(313, 10).Value = 1119
$ws.Cells.Item(313, 11).Value = 1
$ws.Cells.Item(313, 12).Value = 0.6666666666666666
$ws.Cells.Item(313, 14).Value = 1
$ws.Cells.Item(313, 15).Value = "T_1"
$ws.Cells.Item(313, 16).Value = 14
$ws.Cells.Item(313, 17).Value = 3
$ws.Cells.Item(313, 19).Value = 6
$ws.Cells.Item(313, 20).Value = 2
$ws.Cells.Item(313, 22).Value = 28
$ws.Cells.Item(313, 23).Value = 9
$ws.Cells.Item(313, 26).Value = 3
$ws.Cells.Item(313, 27).Value = 3
$ws.Cells.Item(313, 29).Value = 9

# --- Step 4: populate brand-new rows 314-316 ---
# Row 314
$ws.Cells.Item(314, 1).Value = 212
$ws.Cells.Item(314, 2).Value = 2022
$ws.Cells.Item(314, 3).Value = 3
$ws.Cells.Item(314, 4).Value = 31
$ws.Cells.Item(314, 5).Value = 4
$ws.Cells.Item(314, 6).Value = 4
$ws.Cells.Item(314, 7).Value = "ifop"
$ws.Cells.Item(314, 8).Value = "included"
$ws.Cells.Item(314, 9).Value = "rolling"
$ws.Cells.Item(314, 10).Value = 2500
$ws.Cells.Item(314, 11).Value = 1
$ws.Cells.Item(314, 12).Value = 0.6666666666666666
$ws.Cells.Item(314, 13).Value = 0.5
$ws.Cells.Item(314, 14).Value = 1
$ws.Cells.Item(314, 15).Value = 0.5
$ws.Cells.Item(314, 16).Value = 15.5
$ws.Cells.Item(314, 17).Value = 3
$ws.Cells.Item(314, 19).Value = 4.5
$ws.Cells.Item(314, 20).Value = 2
$ws.Cells.Item(314, 22).Value = 27.5
$ws.Cells.Item(314, 23).Value = 10
$ws.Cells.Item(314, 26).Value = 2
$ws.Cells.Item(314, 27).Value = 2
$ws.Cells.Item(314, 28).Value = 22
$ws.Cells.Item(314, 29).Value = 10

# Row 315
$ws.Cells.Item(315, 1).Value = 213
$ws.Cells.Item(315, 2).Value = 2022
$ws.Cells.Item(315, 3).Value = 3
$ws.Cells.Item(315, 4).Value = 31
$ws.Cells.Item(315, 5).Value = 4
$ws.Cells.Item(315, 6).Value = 4
$ws.Cells.Item(315, 7).Value = "ipsos"
$ws.Cells.Item(315, 8).Value = "excluded"
$ws.Cells.Item(315, 9).Value = "rolling"
$ws.Cells.Item(315, 10).Value = 913
$ws.Cells.Item(315, 11).Value = 0
$ws.Cells.Item(315, 12).Value = 1
$ws.Cells.Item(315, 13).Value = 0.5
$ws.Cells.Item(315, 14).Value = 1
$ws.Cells.Item(315, 15).Value = 0.5
$ws.Cells.Item(315, 16).Value = 16
$ws.Cells.Item(315, 17).Value = 3.5
$ws.Cells.Item(315, 19).Value = 5.5
$ws.Cells.Item(315, 20).Value = 2
$ws.Cells.Item(315, 22).Value = 26.5
$ws.Cells.Item(315, 23).Value = 8.5
$ws.Cells.Item(315, 26).Value = 3
$ws.Cells.Item(315, 27).Value = 2
$ws.Cells.Item(315, 28).Value = 21
$ws.Cells.Item(315, 29).Value = 10.5

# Row 316
$ws.Cells.Item(316, 1).Value = 214
$ws.Cells.Item(316, 2).Value = 2022
$ws.Cells.Item(316, 3).Value = 4
$ws.Cells.Item(316, 4).Value = 1
$ws.Cells.Item(316, 5).Value = 4
$ws.Cells.Item(316, 6).Value = 4
$ws.Cells.Item(316, 7).Value = "harris"
$ws.Cells.Item(316, 8).Value = "included"
$ws.Cells.Item(316, 9).Value = "regular"
$ws.Cells.Item(316, 10).Value = 1800
$ws.Cells.Item(316, 11).Value = 1
$ws.Cells.Item(316, 12).Value = 1
$ws.Cells.Item(316, 13).Value = 0.5
$ws.Cells.Item(316, 14).Value = 1
$ws.Cells.Item(316, 15).Value = 0.5
$ws.Cells.Item(316, 16).Value = 17
$ws.Cells.Item(316, 17).Value = 2.5
$ws.Cells.Item(316, 19).Value = 5
$ws.Cells.Item(316, 20).Value = 2
$ws.Cells.Item(316, 22).Value = 26.5
$ws.Cells.Item(316, 23).Value = 9.5
$ws.Cells.Item(316, 26).Value = 2
$ws.Cells.Item(316, 27).Value = 1.5
$ws.Cells.Item(316, 28).Value = 23
$ws.Cells.Item(316, 29).Value = 9.5

# --- Step 5: update view/selection state ---
$ws.Activate()
$ws.Range("Z310").Select()

Write-Host "done"